$wb = $excel.ActiveWorkbook

$wsReview = $wb.Worksheets.Item("LH_TC_NAVIGATION _REVIEWS")

# As the Navigation TC creator, close out the "Owner Status" for the
# reviews already addressed, and mark the one that doesn't apply as
# "NotApplicable". Also fill in the Reviewer/Version for review #003
# (row 4), which had been left blank.
$wsReview.Range("I2").Value = "Closed"
$wsReview.Range("I3").Value = "Closed"
$wsReview.Range("D4").Value = "Ahmed`nAbuzaid"
$wsReview.Range("E4").Value = "v1.0"
$wsReview.Range("I4").Value = "Closed"
$wsReview.Range("I5").Value = "NotApplicable"
$wsReview.Range("I6").Value = "Closed"

# Make the reviews sheet the active tab/view, focused on E4 with the
# zoom level used while wrapping up the review.
$wsReview.Activate()
$excel.ActiveWindow.Zoom = 72
$wsReview.Range("E4").Select() | Out-Null
